# Append: 2025-09-07 06:24 JST
# The scraper re-ran and produced a new snapshot. The sheet keeps only a
# rolling window of rows: row 2 (unchanged job posting, timestamp bumped)
# stays, and three more postings from the new crawl replace the old rows
# 3-10 (several of which - the old rows 5, 6 and 8 - survive with a
# refreshed timestamp; the rest are dropped).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-07 06:24:00"

# --- Row 2: same posting, only the fetch timestamp changes ---------------
$ws.Range("A2").Value = $newTimestamp

# --- Row 3: becomes the old row 5 (Instagram -> Google business MEO tool) -
$ws.Range("A3").Value = $newTimestamp
$ws.Range("B3").Value = "【急募】Instagram投稿を自動でGoogleビジネスに連携するMEOツール"
$ws.Range("D3").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5388589"
$ws.Range("G3").Value = 68
$ws.Range("H3").Value = "◆ツール"

# --- Row 4: becomes the old row 6 (IB報酬 EA開発) -------------------------
$ws.Range("A4").Value = $newTimestamp
$ws.Range("B4").Value = "IB報酬を得るための高性能EA開発依頼"
$ws.Range("D4").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5388547"
$ws.Range("G4").Value = 68
$ws.Range("H4").Value = "◆開発"

# --- Row 5: becomes the old row 8 (限定公開 PR) ---------------------------
$ws.Range("A5").Value = $newTimestamp
$ws.Range("B5").Value = "限定公開 PR 限定公開の仕事"
$ws.Range("D5").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5385681"
$ws.Range("G5").Value = 25
$ws.Range("H5").ClearContents()

# --- Drop the old rows 6-10 (Bubble講師, 運送料金, マネーフォワード, ...) --
$ws.Range("A6:H10").EntireRow.Delete()

# --- Column H got a little narrower now that the remaining skill tags ----
# --- are shorter -----------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 16.166666666666668

# --- Hyperlinks: rebuild F2:F5 so the relationship targets follow the ----
# --- rows they now sit in (engine has no in-place hyperlink-target edit) -
$ws.Range("A1:H5").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5388502")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5388589")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5388547")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5385681")
$ws.Range("F2:F5").Style = "Hyperlink"
